$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locations")
$ws.Range("A3").Value = "BC"
